$wb = $excel.ActiveWorkbook

# The new "Spain" sheet is built from a duplicate of the "Italy" sheet
# (same layout / styles / merged cells), inserted right after it.
$italy = $wb.Worksheets.Item("Italy")
$italy.Copy($null, $italy)

# The freshly inserted copy is now the last worksheet in the workbook.
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# Fill in the market specific values for Spain.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2037"

# Columns were re-fitted to the new (shorter) content.
$spain.Columns("A").ColumnWidth = 25.109375
$spain.Columns("B").ColumnWidth = 15.21875
$spain.Columns("D").ColumnWidth = 21.5546875

# Italy is no longer the active sheet/tab - it keeps a plain A1:D11 selection.
$italy.Activate()
$italy.Range("A1:D11").Select()

# Spain becomes the active sheet, with E5 selected as the last user action.
$spain.Activate()
$spain.Range("E5").Select()
